$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.430.41'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.58%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.114.43'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.17%  '

$ws.Range('E4').Value = '  -0.16%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '346.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.61%  '

$ws.Range('E6').Value = '  -0.09%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5227'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.25%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4459'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.31%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '54.13'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.07%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09418'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.59%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.179'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.13%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.24'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.19%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.186.01'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.36%  '

$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.706'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.55%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.954'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.78%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '101.95'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.83%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001165'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.65%  '

$ws.Range('E18').Value = '  -0.03%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '21.45'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.03%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06730'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.58%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.322'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.05%  '

$ws.Range('E22').Value = '  -0.08%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.497.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.65%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.71'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.16%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.328'
$ws.Range('D25').Style = 'Normal'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.375.05'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.55%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.14'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.17%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.557'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.32%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.37'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.27%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.94'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.87%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.159'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.38%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.777'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.93%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.1060'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.60%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.871'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +12.83%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.283'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.47%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.964'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.18%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.59'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.59%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02646'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.87%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06870'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.31%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.7101'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.68%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.62'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.95%  '

$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.339'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.35%  '

$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.2247'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.55%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6890'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.65%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.60'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.57%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.382'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.21%  '

$ws.Range('E47').Value = '  -0.05%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.397'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +19.59%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.649'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.92%  '

$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00000000350'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.68%  '

$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.227'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.70%  '
